$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3774, 3989, 4136, 4136, 4136, 4220, 4220, 4398, 4557, 4876, 4912, 4912, 5292, 5292)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
